# Automatische test-sync: 2025-07-31 21:46:50
# Adds a new historical response row (row 9) for Testmail #11 "Mijn retour is nog
# steeds niet verwerkt." to the log sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$subjectFull  = "Testmail #11: Mijn retour is nog steeds niet verwerkt."
$replyText    = "Beste klant,`nBedankt voor je bericht. We begrijpen dat het vervelend is dat je retourzending nog niet verwerkt is. Om je verder te kunnen helpen, ontvangen we graag wat aanvullende informatie zoals het ordernummer van de retourzending. Zodra we deze gegevens hebben, zullen we direct voor je aan de slag gaan om het probleem op te lossen.`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$subjectShort = "Mijn retour is nog steeds niet verwerkt."
$sender       = "mailmind.test@zohomail.eu"
$category     = "Retour / Terugbetaling"
$timestamp    = "2025-07-31 21:45:55"
$ja           = "Ja"
$nee          = "Nee"

$row = 9
$ws.Cells.Item($row, 1).Value = $subjectFull
$ws.Cells.Item($row, 2).Value = $replyText
$ws.Cells.Item($row, 3).Value = $subjectShort
$ws.Cells.Item($row, 4).Value = $sender
$ws.Cells.Item($row, 5).Value = $category
$ws.Cells.Item($row, 6).Value = $timestamp
$ws.Cells.Item($row, 7).Value = $ja
$ws.Cells.Item($row, 8).Value = $nee
$ws.Cells.Item($row, 9).Value = $ja
$ws.Cells.Item($row, 10).Value = $nee

# The reply text contains embedded newlines; the COM engine auto-sizes the
# row height as a side effect of setting such a value. Re-running AutoFit
# restores the default row height so the serialized row doesn't carry a
# stray customHeight marker (matching the rest of the sheet's rows).
$ws.Rows.Item($row).EntireRow.AutoFit()
